# Update extrapolation/calibration values for several expiry rows.
# Commit message: "Removing less than USD 5 price from extrapolation
# calibration because it is just a noise" — the underlying model was
# re-run without the noisy sub-$5 price point, producing new values for
# the ABSM1_RN / M1_RN / CM2_RN / CMN3_RN / CMN4_RN columns (D:H) on the
# affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 121946.9108292361
$ws.Range("E5").Value = -0.0082480150923689
$ws.Range("F5").Value = 0.2307169291066302
$ws.Range("G5").Value = -1.015666674325183
$ws.Range("H5").Value = 11.60767055117902

$ws.Range("D7").Value = 123206.4999981723
$ws.Range("E7").Value = -0.02247645405955778
$ws.Range("F7").Value = 0.2631795808157445
$ws.Range("G7").Value = -1.615194553277307
$ws.Range("H7").Value = 15.74634044432639

$ws.Range("D8").Value = 124705.0613436273
$ws.Range("E8").Value = -0.04287349815143904
$ws.Range("F8").Value = 0.2221166801899969
$ws.Range("G8").Value = -0.9041071097604515
$ws.Range("H8").Value = 7.332362305669718

$ws.Range("D9").Value = 126311.4230927655
$ws.Range("E9").Value = -0.07511570758779666
$ws.Range("F9").Value = 0.3479963475284348
$ws.Range("G9").Value = -1.70977567440395
$ws.Range("H9").Value = 10.65412567188488

$ws.Range("D10").Value = 127879.1725034366
$ws.Range("E10").Value = -0.1104730313090675
$ws.Range("F10").Value = 0.443976593596582
$ws.Range("G10").Value = -1.951602191271967
$ws.Range("H10").Value = 9.980662596743176

$ws.Range("D11").Value = 130025.3427526359
$ws.Range("E11").Value = -0.1243643239224585
$ws.Range("F11").Value = 0.4433521965249804
$ws.Range("G11").Value = -1.687984801689959
$ws.Range("H11").Value = 8.061895370126855

$ws.Range("D14").Value = 119651.6177166612
$ws.Range("E14").Value = 0.05023738386003348
$ws.Range("F14").Value = 0.1635972942669265
$ws.Range("G14").Value = -1.789809345182317
$ws.Range("H14").Value = 19.86502897929982

$ws.Range("D16").Value = 119660.6637560399
$ws.Range("E16").Value = 0.04653120129435591
$ws.Range("F16").Value = 0.1633533670516294
$ws.Range("G16").Value = -0.391049677931278
$ws.Range("H16").Value = 8.70028301677544

$ws.Range("D19").Value = 120532.1712405688
$ws.Range("E19").Value = 0.04033182424523526
$ws.Range("F19").Value = 0.1792797382196268
$ws.Range("G19").Value = -0.2539190280963102
$ws.Range("H19").Value = 5.846359812702327
